$d = $word.ActiveDocument

# Helper: insert a new italic paragraph right after the paragraph whose
# text matches $anchorText, containing $newText.
function Add-ItalicParagraphAfter($anchorText, $newText) {
    $found = $d.Content
    $found.Find.ClearFormatting()
    [void]$found.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $found.InsertParagraphAfter()
    $newRng = $found.Next(4, 1)          # wdParagraph = 4, move to next paragraph range
    $newRng.InsertBefore($newText)
    $textOnly = $d.Range($newRng.Start, $newRng.End - 1)
    $textOnly.Font.Italic = $true
}

# 1. Update activation date
[void]$d.Content.Find.Execute("Ativação: 15/07/2016", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# 2. Add English translation after "Objetivos" body paragraph
Add-ItalicParagraphAfter "Estudo de funções especiais em Engenharia Física." "Introduction to complex variable functions and their applications. Present differential equations of interest in physical engineering and develop solution techniques, verifying properties and resolution methods. Study of special functions in Physical Engineering."

# 3. Add English translation after "Programa resumido" body paragraph
Add-ItalicParagraphAfter "Série de Fourier e Transformadas Integrais. Funções especiais." "Functions of a complex variable. Delta function. Partial differential equations in physical engineering: solution methods, solving boundary value problems, applications. Fourier Series and Integral Transforms. Special functions."

# 4. Add English translation after "Programa" body paragraph
Add-ItalicParagraphAfter "Funções especiais: Polinômios de Legendre, Harmônicos Esféricos, Funções de Bessel." "Functions of a complex variable: infinite series, analytical functions, Cauchy Riemann conditions, boundary integrals, Cauchy's theorem, residue theorem, Delta function. Laplace equation, diffusion equation (of heat), wave equation (vibrating string), Fourier series, Fourier and Laplace integral transforms. Special functions: Legendre Polynomials, Spherical Harmonics, Bessel Functions."

# 5. Update grading formula
[void]$d.Content.Find.Execute("Conceito Final = (P1 + 2P2)/3", $true, $false, $false, $false, $false, $true, 1, $false, "Conceito Final = (P1 + P2)/2", 2)
